$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates between the two groups of rows (2-4 and 5-7)
$ws.Range("D2").Value = 44294
$ws.Range("D3").Value = 44294
$ws.Range("D4").Value = 44294
$ws.Range("D5").Value = 44295
$ws.Range("D6").Value = 44295
$ws.Range("D7").Value = 44295

# Swap the Volumen values for M3 and M6
$ws.Range("M3").Value = 240
$ws.Range("M6").Value = 200
